$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row before row 27 (pushes the former row 27 and
# everything below it down by one row). This also extends the
# C14:C29 / D14:D29 merged block to C14:C30 / D14:D30, since row 27
# falls inside that merged range.
$ws.Rows("27:27").Insert(-4121)   # xlShiftDown = -4121

# The insert operation leaves the new row 27 mostly formatted like the
# row that used to be there (style carries over), but the "-" marker in
# column F gets dropped - restore it explicitly.
$ws.Range("F27").Value = "-"

# Populate the newly inserted row with the new sub-task description.
$ws.Range("E27").Value = 'Implémenter le bouton "Modifier un indice"'

# The conditional formatting range (originally G71:G1048576) needs to
# shift down by one row as well, since it refers to rows below the
# insertion point.
$fc = $ws.Cells.FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("G72:G1048576"))
$fc.Formula1 = '=NOT(ISERROR(SEARCH("Oui",G72)))'

# Restore selection to match the saved state of the workbook.
$ws.Range("F28").Select()
